# Insert a new "Assurance MM1 DCA" row above the existing "Project MM18" row
# (currently row 19), shifting the Project MM18 block down by one row, then
# populate the new row with the colour values and fix up the saved
# selection to match the committed workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 19; existing rows 19 and below move down to 20+.
# Using EntireRow.Insert copies the formatting down from the row above,
# which reproduces the D-column "highlight" style used throughout this
# table.
$ws.Rows.Item(19).Insert()

$ws.Range("A19").Value = "Assurance MM1 DCA"
$ws.Range("B19").Value = "Bright Yellow"
$ws.Range("C19").Value = "Magnolia"
$ws.Range("D19").Value = "Duck Egg"
$ws.Range("E19").Value = "Salmon Pink"
$ws.Range("F19").Value = "White "
$ws.Range("G19").Value = "Purple"

# Restore the expected active selection on the sheet.
$ws.Range("K25").Select()
